# Add a new "PUBLONS" row to the Test Suite worksheet, mirroring the
# existing rows (TSID / Description / Runmode columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

$newRow = 19

$ws.Range("A$newRow").Value = "PUBLONS"
$ws.Range("B$newRow").Value = "PUBLONS module"
$ws.Range("C$newRow").Value = "Y"

# Match the formatting used by the previous last row (thin black box
# border around each cell, no fill) so the new row looks consistent
# with the rest of the table.
$newRange = $ws.Range("A$($newRow):C$newRow")
$newRange.Borders.Color = 0
$newRange.Borders.LineStyle = 1

# Leave the new row's last cell selected, like the source workbook does
# after appending a row.
$ws.Range("B$newRow").Select()
